# Applies the two substantive edits captured by the commit:
#  1. The table on slide 6 ("SOURCES OF FINANCE") is switched from the
#     deck's custom table style to the built-in PowerPoint table style
#     {5F335975-6203-4C09-AD4C-2E1A0C02A808}.
#  2. The presentation's theme (Design) colour scheme is changed from the
#     "Integral" palette to the stock "Office Theme" palette - i.e. the
#     same effect as picking the default "Office Theme" entry in the
#     PowerPoint Design gallery.

$p = $ppt.ActivePresentation

# --- 1. Re-style the sources-of-finance table on slide 6 -------------------
$slide = $p.Slides.Item(6)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{5F335975-6203-4C09-AD4C-2E1A0C02A808}")
    }
}

# --- 2. Swap the theme colour scheme from "Integral" to "Office Theme" -----
# Theme.ThemeColorScheme order is:
#  1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3, 8 accent4,
#  9 accent5, 10 accent6, 11 hlink, 12 folHlink
# .RGB uses the standard OLE BGR-packed integer (R + G*256 + B*65536).
$officeThemeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Item($i).RGB = $officeThemeColors[$i - 1]
}
